# Case_4_85 / res_bus / vm_pu.xlsx: update bus voltage-magnitude results
# for the 380 kV slack-bus case (slack vm_pu 1.05 -> 1.02, rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.045602806638494
$ws.Cells.Item(2, 4).Value = 1.043903024682111
$ws.Cells.Item(2, 5).Value = 1.058916436070497
$ws.Cells.Item(2, 6).Value = 1.065712182328834
$ws.Cells.Item(2, 9).Value = 1.037238286325926
$ws.Cells.Item(2, 10).Value = 1.050661661101414
$ws.Cells.Item(2, 11).Value = 1.046675580715689
$ws.Cells.Item(2, 12).Value = 1.061647363672588
$ws.Cells.Item(2, 13).Value = 1.06842468439727
$ws.Cells.Item(2, 14).Value = 1.020625096886802
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046930593045064
$ws.Cells.Item(3, 4).Value = 1.04485803872579
$ws.Cells.Item(3, 5).Value = 1.06022426755055
$ws.Cells.Item(3, 6).Value = 1.06707252464011
$ws.Cells.Item(3, 9).Value = 1.037504011927837
$ws.Cells.Item(3, 10).Value = 1.05163530968464
$ws.Cells.Item(3, 11).Value = 1.047441348312247
$ws.Cells.Item(3, 12).Value = 1.062768072140542
$ws.Cells.Item(3, 13).Value = 1.069599117913084
$ws.Cells.Item(3, 14).Value = 1.020958805714661
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.047789145348526
$ws.Cells.Item(4, 4).Value = 1.045475257047267
$ws.Cells.Item(4, 5).Value = 1.061070261394275
$ws.Cells.Item(4, 6).Value = 1.06795252108649
$ws.Cells.Item(4, 9).Value = 1.037674299964221
$ws.Cells.Item(4, 10).Value = 1.052264245471713
$ws.Cells.Item(4, 11).Value = 1.047935493610569
$ws.Cells.Item(4, 12).Value = 1.063492436559752
$ws.Cells.Item(4, 13).Value = 1.070358278964812
$ws.Cells.Item(4, 14).Value = 1.021174118593319
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.048149937557373
$ws.Cells.Item(5, 4).Value = 1.045734560394307
$ws.Cells.Item(5, 5).Value = 1.061425858059503
$ws.Cells.Item(5, 6).Value = 1.068322418516485
$ws.Cells.Item(5, 9).Value = 1.037745493764397
$ws.Cells.Item(5, 10).Value = 1.052528394096019
$ws.Cells.Item(5, 11).Value = 1.048142909130426
$ws.Cells.Item(5, 12).Value = 1.063796768601885
$ws.Cells.Item(5, 13).Value = 1.070677247538056
$ws.Cells.Item(5, 14).Value = 1.021264488671398
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.04821050787333
$ws.Cells.Item(6, 4).Value = 1.045778088344446
$ws.Cells.Item(6, 5).Value = 1.061485560929766
$ws.Cells.Item(6, 6).Value = 1.068384522909446
$ws.Cells.Item(6, 9).Value = 1.037757424357775
$ws.Cells.Item(6, 10).Value = 1.052572730875818
$ws.Cells.Item(6, 11).Value = 1.048177716211575
$ws.Cells.Item(6, 12).Value = 1.063847856188527
$ws.Cells.Item(6, 13).Value = 1.070730793137506
$ws.Cells.Item(6, 14).Value = 1.021279653579836
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.047793966829444
$ws.Cells.Item(7, 4).Value = 1.045478722557569
$ws.Cells.Item(7, 5).Value = 1.0610750131227
$ws.Cells.Item(7, 6).Value = 1.067957463878349
$ws.Cells.Item(7, 9).Value = 1.037675252811585
$ws.Cells.Item(7, 10).Value = 1.052267776041209
$ws.Cells.Item(7, 11).Value = 1.047938266374439
$ws.Cells.Item(7, 12).Value = 1.063496503805356
$ws.Cells.Item(7, 13).Value = 1.070362541752745
$ws.Cells.Item(7, 14).Value = 1.021175326701203
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.046051667803576
$ws.Cells.Item(8, 4).Value = 1.044225930399971
$ws.Cells.Item(8, 5).Value = 1.059358479747226
$ws.Cells.Item(8, 6).Value = 1.06617196730815
$ws.Cells.Item(8, 9).Value = 1.037328432488898
$ws.Cells.Item(8, 10).Value = 1.050990935355115
$ws.Cells.Item(8, 11).Value = 1.046934657655495
$ws.Cells.Item(8, 12).Value = 1.062026281578409
$ws.Cells.Item(8, 13).Value = 1.068821752201031
$ws.Cells.Item(8, 14).Value = 1.020738003975323
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.04297663538363
$ws.Cells.Item(9, 4).Value = 1.042012603148095
$ws.Cells.Item(9, 5).Value = 1.056331587529359
$ws.Cells.Item(9, 6).Value = 1.063023733133128
$ws.Cells.Item(9, 9).Value = 1.036704588685168
$ws.Cells.Item(9, 10).Value = 1.048732592550515
$ws.Cells.Item(9, 11).Value = 1.045155690988951
$ws.Cells.Item(9, 12).Value = 1.059429236947078
$ws.Cells.Item(9, 13).Value = 1.066100604792874
$ws.Cells.Item(9, 14).Value = 1.019962610639467
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.04092308796209
$ws.Cells.Item(10, 4).Value = 1.040533069110926
$ws.Cells.Item(10, 5).Value = 1.05431201427735
$ws.Cells.Item(10, 6).Value = 1.060923374695313
$ws.Cells.Item(10, 9).Value = 1.036280105414806
$ws.Cells.Item(10, 10).Value = 1.047221228980953
$ws.Cells.Item(10, 11).Value = 1.043962549364168
$ws.Cells.Item(10, 12).Value = 1.057693437184946
$ws.Cells.Item(10, 13).Value = 1.064282226203721
$ws.Cells.Item(10, 14).Value = 1.019442423362833
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.040032987766163
$ws.Cells.Item(11, 4).Value = 1.0398914428082
$ws.Cells.Item(11, 5).Value = 1.053437080126296
$ws.Cells.Item(11, 6).Value = 1.060013483052992
$ws.Cells.Item(11, 9).Value = 1.036094251750476
$ws.Cells.Item(11, 10).Value = 1.046565378834663
$ws.Cells.Item(11, 11).Value = 1.043444179834364
$ws.Cells.Item(11, 12).Value = 1.056940722801383
$ws.Cells.Item(11, 13).Value = 1.063493789357962
$ws.Cells.Item(11, 14).Value = 1.019216392670412
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.039702224773093
$ws.Cells.Item(12, 4).Value = 1.039652964816869
$ws.Cells.Item(12, 5).Value = 1.05311201958805
$ws.Cells.Item(12, 6).Value = 1.059675440983058
$ws.Cells.Item(12, 9).Value = 1.036024908492312
$ws.Cells.Item(12, 10).Value = 1.046321550315308
$ws.Cells.Item(12, 11).Value = 1.043251371792438
$ws.Cells.Item(12, 12).Value = 1.056660961981711
$ws.Cells.Item(12, 13).Value = 1.063200764666262
$ws.Cells.Item(12, 14).Value = 1.019132315703141
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.039773180977916
$ws.Cells.Item(13, 4).Value = 1.039704125956133
$ws.Cells.Item(13, 5).Value = 1.053181749457485
$ws.Cells.Item(13, 6).Value = 1.059747955292203
$ws.Cells.Item(13, 9).Value = 1.036039796852561
$ws.Cells.Item(13, 10).Value = 1.046373862213534
$ws.Cells.Item(13, 11).Value = 1.043292741690208
$ws.Cells.Item(13, 12).Value = 1.056720979342307
$ws.Cells.Item(13, 13).Value = 1.063263626947127
$ws.Cells.Item(13, 14).Value = 1.019150355903726
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.04000564968105
$ws.Cells.Item(14, 4).Value = 1.039871733205194
$ws.Cells.Item(14, 5).Value = 1.053410211993696
$ws.Cells.Item(14, 6).Value = 1.059985541798186
$ws.Cells.Item(14, 9).Value = 1.036088526123063
$ws.Cells.Item(14, 10).Value = 1.046545228330211
$ws.Cells.Item(14, 11).Value = 1.043428247630728
$ws.Cells.Item(14, 12).Value = 1.056917601157366
$ws.Cells.Item(14, 13).Value = 1.063469571216324
$ws.Cells.Item(14, 14).Value = 1.019209445277515
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.040148862605211
$ws.Cells.Item(15, 4).Value = 1.039974981691112
$ws.Cells.Item(15, 5).Value = 1.053550965766497
$ws.Cells.Item(15, 6).Value = 1.060131917572168
$ws.Cells.Item(15, 9).Value = 1.03611850886831
$ws.Cells.Item(15, 10).Value = 1.046650783854238
$ws.Cells.Item(15, 11).Value = 1.043511702454539
$ws.Cells.Item(15, 12).Value = 1.057038723818302
$ws.Cells.Item(15, 13).Value = 1.06359643841168
$ws.Cells.Item(15, 14).Value = 1.019245836378687
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.040982141593006
$ws.Cells.Item(16, 4).Value = 1.040575630916333
$ws.Cells.Item(16, 5).Value = 1.054370071038332
$ws.Cells.Item(16, 6).Value = 1.060983751936038
$ws.Cells.Item(16, 9).Value = 1.036292396642184
$ws.Cells.Item(16, 10).Value = 1.047264725403822
$ws.Cells.Item(16, 11).Value = 1.043996915160527
$ws.Cells.Item(16, 12).Value = 1.057743368836194
$ws.Cells.Item(16, 13).Value = 1.064334529353654
$ws.Cells.Item(16, 14).Value = 1.019457407627363
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.041504590956122
$ws.Cells.Item(17, 4).Value = 1.040952138672981
$ws.Cells.Item(17, 5).Value = 1.054883751951528
$ws.Cells.Item(17, 6).Value = 1.061517968952093
$ws.Cells.Item(17, 9).Value = 1.036400922405208
$ws.Cells.Item(17, 10).Value = 1.047649452305986
$ws.Cells.Item(17, 11).Value = 1.044300811051347
$ws.Cells.Item(17, 12).Value = 1.058185076512134
$ws.Cells.Item(17, 13).Value = 1.064797225945223
$ws.Cells.Item(17, 14).Value = 1.019589909504591
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.041809240263717
$ws.Cells.Item(18, 4).Value = 1.041171655051943
$ws.Cells.Item(18, 5).Value = 1.055183330621341
$ws.Cells.Item(18, 6).Value = 1.061829528220141
$ws.Cells.Item(18, 9).Value = 1.036464025957736
$ws.Cells.Item(18, 10).Value = 1.047873720236742
$ws.Cells.Item(18, 11).Value = 1.044477901498483
$ws.Cells.Item(18, 12).Value = 1.05844261108507
$ws.Cells.Item(18, 13).Value = 1.06506700635887
$ws.Cells.Item(18, 14).Value = 1.019667119821604
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.041913103267285
$ws.Cells.Item(19, 4).Value = 1.041246488555102
$ws.Cells.Item(19, 5).Value = 1.055285472083573
$ws.Cells.Item(19, 6).Value = 1.061935755142892
$ws.Cells.Item(19, 9).Value = 1.036485509136972
$ws.Cells.Item(19, 10).Value = 1.04795016670666
$ws.Cells.Item(19, 11).Value = 1.044538256530233
$ws.Cells.Item(19, 12).Value = 1.05853040584321
$ws.Cells.Item(19, 13).Value = 1.065158977207629
$ws.Cells.Item(19, 14).Value = 1.019693433719485
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.041448546074378
$ws.Cells.Item(20, 4).Value = 1.040911752718737
$ws.Cells.Item(20, 5).Value = 1.054828643300255
$ws.Cells.Item(20, 6).Value = 1.061460656741632
$ws.Cells.Item(20, 9).Value = 1.036389299073486
$ws.Cells.Item(20, 10).Value = 1.047608188921151
$ws.Cells.Item(20, 11).Value = 1.044268223178258
$ws.Cells.Item(20, 12).Value = 1.058137696454967
$ws.Cells.Item(20, 13).Value = 1.064747593629685
$ws.Cells.Item(20, 14).Value = 1.019575701153634
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.039937197349893
$ws.Cells.Item(21, 4).Value = 1.039822381192052
$ws.Cells.Item(21, 5).Value = 1.05334293746502
$ws.Cells.Item(21, 6).Value = 1.059915580400253
$ws.Cells.Item(21, 9).Value = 1.03607418510356
$ws.Cells.Item(21, 10).Value = 1.046494771284889
$ws.Cells.Item(21, 11).Value = 1.043388351774817
$ws.Cells.Item(21, 12).Value = 1.056859705633393
$ws.Cells.Item(21, 13).Value = 1.063408930324212
$ws.Cells.Item(21, 14).Value = 1.019192048225549
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.038986139151185
$ws.Cells.Item(22, 4).Value = 1.039136584393912
$ws.Cells.Item(22, 5).Value = 1.052408402193446
$ws.Cells.Item(22, 6).Value = 1.058943735153016
$ws.Cells.Item(22, 9).Value = 1.035874272557332
$ws.Cells.Item(22, 10).Value = 1.045793466766175
$ws.Cells.Item(22, 11).Value = 1.042833621650681
$ws.Cells.Item(22, 12).Value = 1.056055201406084
$ws.Cells.Item(22, 13).Value = 1.062566307976876
$ws.Cells.Item(22, 14).Value = 1.018950140807183
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.039490391807189
$ws.Cells.Item(23, 4).Value = 1.039500221125972
$ws.Cells.Item(23, 5).Value = 1.052903857545349
$ws.Cells.Item(23, 6).Value = 1.059458967538603
$ws.Cells.Item(23, 9).Value = 1.035980419787401
$ws.Cells.Item(23, 10).Value = 1.046165361584345
$ws.Cells.Item(23, 11).Value = 1.043127839474888
$ws.Cells.Item(23, 12).Value = 1.056481778506653
$ws.Cells.Item(23, 13).Value = 1.063013089516157
$ws.Cells.Item(23, 14).Value = 1.019078446204841
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.041473870611247
$ws.Cells.Item(24, 4).Value = 1.040930001684164
$ws.Cells.Item(24, 5).Value = 1.054853544659463
$ws.Cells.Item(24, 6).Value = 1.06148655378746
$ws.Cells.Item(24, 9).Value = 1.036394551767574
$ws.Cells.Item(24, 10).Value = 1.047626834491033
$ws.Cells.Item(24, 11).Value = 1.044282948750946
$ws.Cells.Item(24, 12).Value = 1.058159105790848
$ws.Cells.Item(24, 13).Value = 1.064770020653782
$ws.Cells.Item(24, 14).Value = 1.019582121529983
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.043772208361941
$ws.Cells.Item(25, 4).Value = 1.042585494518221
$ws.Cells.Item(25, 5).Value = 1.057114386554612
$ws.Cells.Item(25, 6).Value = 1.063837881821799
$ws.Cells.Item(25, 9).Value = 1.036867377366674
$ws.Cells.Item(25, 10).Value = 1.049317438873305
$ws.Cells.Item(25, 11).Value = 1.045616850511063
$ws.Cells.Item(25, 12).Value = 1.060101403688621
$ws.Cells.Item(25, 13).Value = 1.066804826860655
$ws.Cells.Item(25, 14).Value = 1.020163639169905
